$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "315.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.65%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.97%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.156"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.92%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08157"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.69%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.130"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.36%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9294"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.23%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1014"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.99%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1874"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.02%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09218"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "7.34%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03617"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.82%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09930"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.08%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001435"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.46%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005693"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.59%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.455"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.09%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.153"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.27%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "8.85%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1337"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.00%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.183"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.54%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2194"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.42%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04612"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.38%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001252"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.39%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004726"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.65%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001253"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-21.84%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004525"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-4.82%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01974"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.24%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04922"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.16%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007892"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.51%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1401"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.22%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007853"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.51%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002117"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.61%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01184"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.54%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006544"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.47%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.45%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.03"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-25.98%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001910"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.56%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002112"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.45%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002011"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.45%"
